$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; this shifts every existing row down by one.
$ws.Rows.Item(1).Insert()

# New header row (row 1): tracking columns for model-making progress.
$ws.Range("D1").Value = "Placeholder"
$ws.Range("E1").Value = "Final"

# Mark Chair / Bed / Table / Coffin (now rows 3-6) as placeholder models.
$ws.Range("D3").Value = "x"

$ws.Range("F1").Value = "Needs more coding"

$ws.Range("D4").Value = "x"
$ws.Range("D5").Value = "x"
$ws.Range("D6").Value = "x"

# Mark FarmPlot (now row 7) as needing more coding.
$ws.Range("F7").Value = "x"

# Auto-fit the new columns to their content.
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(6).AutoFit()

# Restore the selection to where the author left off.
$ws.Range("F8").Select()
